# Auto-generated Excel COM-interop script applying the cryptos.xlsx data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Sheet, [string]$Address, [string]$Text)
    # Force literal text even for values that look numeric (e.g. "2.55", "9.84")
    # by using a leading apostrophe, then reset the cell style to "Normal" so no
    # stray number-format / quote-prefix style sticks to the cell (matches source).
    $Sheet.Range($Address).Value = "'" + $Text
    $Sheet.Range($Address).Style = "Normal"
}

# Row 2
Set-TextCell $ws "D2" "43.731.38"
Set-TextCell $ws "E2" "  +4.48%  "

# Row 3
Set-TextCell $ws "D3" "2.265.96"
Set-TextCell $ws "E3" "  +1.57%  "

# Row 4
Set-TextCell $ws "E4" "  -0.14%  "

# Row 5
Set-TextCell $ws "D5" "230.73"
Set-TextCell $ws "E5" "  -0.32%  "

# Row 6
Set-TextCell $ws "E6" "  +0.04%  "

# Row 7
Set-TextCell $ws "D7" "61.11"
Set-TextCell $ws "E7" "  -1.12%  "

# Row 8
Set-TextCell $ws "E8" "  -0.13%  "

# Row 9
Set-TextCell $ws "E9" "  +4.65%  "

# Row 10
Set-TextCell $ws "D10" "58.06"
Set-TextCell $ws "E10" "  -2.17%  "

# Row 11
Set-TextCell $ws "D11" "0.0933"
Set-TextCell $ws "E11" "  +4.55%  "

# Row 12
Set-TextCell $ws "E12" "  +0.80%  "

# Row 13
Set-TextCell $ws "D13" "2.604.00"
Set-TextCell $ws "E13" "  +1.44%  "

# Row 14
Set-TextCell $ws "D14" "15.58"
Set-TextCell $ws "E14" "  -0.37%  "

# Row 15
Set-TextCell $ws "D15" "23.51"
Set-TextCell $ws "E15" "  +6.96%  "

# Row 16
Set-TextCell $ws "E16" "  +3.77%  "

# Row 17
Set-TextCell $ws "E17" "  +1.08%  "

# Row 18
Set-TextCell $ws "D18" "2.263.67"
Set-TextCell $ws "E18" "  +0.31%  "

# Row 19
Set-TextCell $ws "D19" "42.860.09"
Set-TextCell $ws "E19" "  +2.62%  "

# Row 20
Set-TextCell $ws "D20" "0.0₃0934"
Set-TextCell $ws "E20" "  +4.02%  "

# Row 21
Set-TextCell $ws "D21" "72.86"
Set-TextCell $ws "E21" "  +0.82%  "

# Row 22
Set-TextCell $ws "E22" "  +2.53%  "

# Row 23
Set-TextCell $ws "D23" "253.56"
Set-TextCell $ws "E23" "  +1.52%  "

# Row 24
Set-TextCell $ws "E24" "  -0.08%  "

# Row 25
Set-TextCell $ws "D25" "2.55"
Set-TextCell $ws "E25" "  +6.32%  "

# Row 26
Set-TextCell $ws "D26" "2.29"
Set-TextCell $ws "E26" "  -3.02%  "

# Row 27
Set-TextCell $ws "D27" "9.84"
Set-TextCell $ws "E27" "  +1.59%  "

# Row 28
Set-TextCell $ws "D28" "170.52"
Set-TextCell $ws "E28" "  +2.23%  "

# Row 29
Set-TextCell $ws "E29" "  -1.53%  "

# Row 30
Set-TextCell $ws "D30" "20.46"
Set-TextCell $ws "E30" "  +2.53%  "

# Row 31
Set-TextCell $ws "E31" "  +1.54%  "

# Row 32
Set-TextCell $ws "D32" "2.66"
Set-TextCell $ws "E32" "  +0.72%  "

# Row 33
Set-TextCell $ws "E33" "  -0.05%  "

# Row 34
Set-TextCell $ws "D34" "5.05"
Set-TextCell $ws "E34" "  +0.99%  "

# Row 35
Set-TextCell $ws "D35" "4.78"
Set-TextCell $ws "E35" "  +2.12%  "

# Row 36
Set-TextCell $ws "E36" "  +3.71%  "

# Row 37
Set-TextCell $ws "D37" "6.46"
Set-TextCell $ws "E37" "  -2.66%  "

# Row 38
Set-TextCell $ws "E38" "  +1.43%  "

# Row 39
Set-TextCell $ws "D39" "3.59"
Set-TextCell $ws "E39" "  -1.39%  "

# Row 40
Set-TextCell $ws "E40" "  +4.14%  "

# Row 41
Set-TextCell $ws "E41" "  +0.14%  "

# Row 42
Set-TextCell $ws "B42" "FraxShare"
Set-TextCell $ws "C42" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell $ws "D42" "8.71"
Set-TextCell $ws "E42" "  +1.67%  "

# Row 43
Set-TextCell $ws "B43" "TerraClassic"
Set-TextCell $ws "C43" "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
Set-TextCell $ws "D43" "0.000227"
Set-TextCell $ws "E43" "  -10.67%  "

# Row 44
Set-TextCell $ws "D44" "0.0989"
Set-TextCell $ws "E44" "  +1.04%  "

# Row 45
Set-TextCell $ws "D45" "4.52"
Set-TextCell $ws "E45" "  -6.35%  "

# Row 46
Set-TextCell $ws "E46" "  -0.76%  "

# Row 47
Set-TextCell $ws "D47" "98.21"
Set-TextCell $ws "E47" "  -0.66%  "

# Row 48
Set-TextCell $ws "D48" "1.471.43"
Set-TextCell $ws "E48" "  -0.54%  "

# Row 49
Set-TextCell $ws "D49" "16.63"
Set-TextCell $ws "E49" "  +0.97%  "

# Row 50
Set-TextCell $ws "E50" "  +0.67%  "

# Row 51
Set-TextCell $ws "D51" "2.25"
Set-TextCell $ws "E51" "  +7.36%  "
